$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.19%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '34.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '12.65%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.10%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07827'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '6.19%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.319'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.23%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.067'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.63%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.989'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '6.76%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9256'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.53%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1009'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9.02%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1826'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '7.63%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.86%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03422'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '9.97%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09909'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.55%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001482'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.87%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005759'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.09%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.475'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.15%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.102'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.98%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3418'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.89%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1326'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.89%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.544'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '9.03%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2271'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '6.96%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04650'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.08%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.47%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004338'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.69%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.17%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003399'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.18%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01754'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '12.49%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04750'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '6.11%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007777'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '5.52%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '6.19%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-22.27%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002291'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.95%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009964'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '10.95%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006075'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.53%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.06%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.880'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '48.88%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '34.53%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.06%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.06%'
